$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.064679026603699
$ws.Range("B1").Value = 2.420620203018188
$ws.Range("C1").Value = 5.146184921264648
$ws.Range("D1").Value = 2.30472731590271
$ws.Range("E1").Value = 1.307539224624634
